$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 41.272728
$ws.Range("I8").Value = 41.272728
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 123.818184
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 15.181816
$ws.Range("N8").ClearContents()
$ws.Range("H64").Value = 4977.778
$ws.Range("I64").Value = 4600
$ws.Range("J64").Value = 5166.6665
$ws.Range("K64").Value = 4600
$ws.Range("L64").Value = 5166.6665
$ws.Range("M64").Value = -4352
$ws.Range("N64").Value = -5662.6665
$ws.Range("H67").Value = 4977.778
$ws.Range("I67").Value = 4600
$ws.Range("J67").Value = 5166.6665
$ws.Range("K67").Value = 4600
$ws.Range("L67").Value = 5166.6665
$ws.Range("M67").Value = -3742
$ws.Range("N67").Value = -6882.6665
$ws.Range("H74").Value = 4483.3335
$ws.Range("I74").Value = 3900
$ws.Range("K74").Value = 3900
$ws.Range("M74").Value = -2964
$ws.Range("H76").Value = 7805830.5
$ws.Range("J76").Value = 14224973
$ws.Range("L76").Value = 14224973
$ws.Range("N76").Value = -14225603
$ws.Range("H77").Value = 4483.3335
$ws.Range("I77").Value = 3900
$ws.Range("K77").Value = 19500
$ws.Range("M77").Value = -14820
$ws.Range("H79").Value = 7805830.5
$ws.Range("J79").Value = 14224973
$ws.Range("L79").Value = 14224973
$ws.Range("N79").Value = -14227157
$ws.Range("H107").Value = 317
$ws.Range("J107").Value = 113
$ws.Range("L107").Value = 113
$ws.Range("N107").Value = -3953
$ws.Range("H132").Value = 3231.3022
$ws.Range("I132").Value = 2946
$ws.Range("K132").Value = 8838
$ws.Range("M132").Value = -6308

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1339.4445
$ws.Range("I2").Value = 1380.8667
$ws.Range("K2").Value = 1380.8667
$ws.Range("M2").Value = -1267.8667
$ws.Range("H74").Value = 1840.409
$ws.Range("I74").Value = 1624.8379
$ws.Range("J74").Value = 2979.8572
$ws.Range("K74").Value = 1624.8379
$ws.Range("L74").Value = 2979.8572
$ws.Range("M74").Value = -750.8379
$ws.Range("N74").Value = -4727.8572
$ws.Range("H77").Value = 1840.409
$ws.Range("I77").Value = 1624.8379
$ws.Range("J77").Value = 2979.8572
$ws.Range("K77").Value = 8124.1895
$ws.Range("L77").Value = 14899.286
$ws.Range("M77").Value = -3756.1895
$ws.Range("N77").Value = -23635.286
$ws.Range("H97").Value = 648.6
$ws.Range("I97").Value = 642.125
$ws.Range("K97").Value = 642.125
$ws.Range("M97").Value = -146.125
$ws.Range("H116").Value = 1339.4445
$ws.Range("I116").Value = 1380.8667
$ws.Range("K116").Value = 1380.8667
$ws.Range("M116").Value = 913.1333
$ws.Range("H132").Value = 55558104
$ws.Range("I132").Value = 83335800
$ws.Range("K132").Value = 250007400
$ws.Range("M132").Value = -250004870

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1339.4445
$ws.Range("I3").Value = 1380.8667
$ws.Range("K3").Value = 1380.8667
$ws.Range("M3").Value = -1266.8667
$ws.Range("H81").Value = 82220.89
$ws.Range("J81").Value = 82220.89
$ws.Range("L81").Value = 82220.89
$ws.Range("N81").Value = -84342.89
$ws.Range("H84").Value = 82220.89
$ws.Range("J84").Value = 82220.89
$ws.Range("L84").Value = 246662.67
$ws.Range("N84").Value = -257270.67
$ws.Range("H86").Value = 5789.909
$ws.Range("I86").Value = 6132.778
$ws.Range("K86").Value = 6132.778
$ws.Range("M86").Value = -5009.778
$ws.Range("H89").Value = 5789.909
$ws.Range("I89").Value = 6132.778
$ws.Range("K89").Value = 30663.89
$ws.Range("M89").Value = -25047.89
$ws.Range("H105").Value = 1703.5555
$ws.Range("I105").Value = 1616.5
$ws.Range("K105").Value = 1616.5
$ws.Range("M105").Value = 130.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 2636.25
$ws.Range("I39").Value = 2636.25
$ws.Range("K39").Value = 2636.25
$ws.Range("M39").Value = -2245.25
$ws.Range("H49").Value = 2636.25
$ws.Range("I49").Value = 2636.25
$ws.Range("K49").Value = 2636.25
$ws.Range("M49").Value = -2454.25
$ws.Range("H62").Value = 3261.4546
$ws.Range("J62").Value = 3219.5
$ws.Range("L62").Value = 3219.5
$ws.Range("N62").Value = -4467.5
$ws.Range("H65").Value = 3261.4546
$ws.Range("J65").Value = 3219.5
$ws.Range("L65").Value = 16097.5
$ws.Range("N65").Value = -22337.5
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").ClearContents()
$ws.Range("H99").Value = 2592.5293
$ws.Range("J99").Value = 3326
$ws.Range("L99").Value = 3326
$ws.Range("N99").Value = -6322
$ws.Range("H126").Value = 2592.5293
$ws.Range("J126").Value = 3326
$ws.Range("L126").Value = 9978
$ws.Range("N126").Value = -14918

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1064.5714
$ws.Range("I5").Value = 1190.6
$ws.Range("J5").Value = 749.5
$ws.Range("K5").Value = 3571.8
$ws.Range("L5").Value = 2248.5
$ws.Range("M5").Value = -3459.8
$ws.Range("N5").Value = -2472.5
$ws.Range("H17").Value = 498
$ws.Range("I17").Value = 529.4
$ws.Range("J17").Value = 445.66666
$ws.Range("K17").Value = 1588.2
$ws.Range("L17").Value = 1336.99998
$ws.Range("M17").Value = -1419.2
$ws.Range("N17").Value = -1674.99998
$ws.Range("H34").Value = 1691.3684
$ws.Range("I34").Value = 216.33333
$ws.Range("J34").Value = 1967.9375
$ws.Range("K34").Value = 648.99999
$ws.Range("L34").Value = 5903.8125
$ws.Range("M34").Value = -564.99999
$ws.Range("N34").Value = -6071.8125
$ws.Range("H39").Value = 2829.55
$ws.Range("I39").Value = 454.77777
$ws.Range("K39").Value = 1364.33331
$ws.Range("M39").Value = -1070.33331
$ws.Range("H55").Value = 538.04
$ws.Range("I55").Value = 434.13635
$ws.Range("J55").Value = 1300
$ws.Range("K55").Value = 1302.40905
$ws.Range("L55").Value = 3900
$ws.Range("M55").Value = -1125.40905
$ws.Range("N55").Value = -4254
$ws.Range("H70").Value = 3870.6
$ws.Range("I70").Value = 3870.6
$ws.Range("K70").Value = 11611.8
$ws.Range("M70").Value = -11296.8
$ws.Range("H73").Value = 3870.6
$ws.Range("I73").Value = 3870.6
$ws.Range("K73").Value = 11611.8
$ws.Range("M73").Value = -10519.8
$ws.Range("H135").Value = 1064.5714
$ws.Range("I135").Value = 1190.6
$ws.Range("J135").Value = 749.5
$ws.Range("K135").Value = 10715.4
$ws.Range("L135").Value = 6745.5
$ws.Range("M135").Value = -8180.4
$ws.Range("N135").Value = -11815.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7789.84
$ws.Range("I70").Value = 7355.4736
$ws.Range("K70").Value = 7355.4736
$ws.Range("M70").Value = -7085.4736
$ws.Range("H73").Value = 7789.84
$ws.Range("I73").Value = 7355.4736
$ws.Range("K73").Value = 7355.4736
$ws.Range("M73").Value = -6419.4736
$ws.Range("H80").Value = 2810.25
$ws.Range("I80").Value = 2848.3635
$ws.Range("J80").Value = 2726.4
$ws.Range("K80").Value = 2848.3635
$ws.Range("L80").Value = 2726.4
$ws.Range("M80").Value = -1850.3635
$ws.Range("N80").Value = -4722.4
$ws.Range("H83").Value = 2810.25
$ws.Range("I83").Value = 2848.3635
$ws.Range("J83").Value = 2726.4
$ws.Range("K83").Value = 14241.8175
$ws.Range("L83").Value = 13632
$ws.Range("M83").Value = -9249.817499999999
$ws.Range("N83").Value = -23616
$ws.Range("H97").Value = 962.7895
$ws.Range("I97").Value = 734.4286
$ws.Range("J97").Value = 1602.2
$ws.Range("K97").Value = 734.4286
$ws.Range("L97").Value = 1602.2
$ws.Range("M97").Value = -238.4286
$ws.Range("N97").Value = -2594.2
$ws.Range("H123").Value = 100000
$ws.Range("J123").Value = 100000
$ws.Range("L123").Value = 100000
$ws.Range("N123").Value = -104900
$ws.Range("H132").Value = 3763.3845
$ws.Range("I132").Value = 3299.125
$ws.Range("K132").Value = 9897.375
$ws.Range("M132").Value = -7367.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 19998.666
$ws.Range("J39").Value = 19997
$ws.Range("L39").Value = 19997
$ws.Range("N39").Value = -20917
$ws.Range("H46").Value = 2660.5
$ws.Range("I46").Value = 1341.3334
$ws.Range("J46").Value = 3649.875
$ws.Range("K46").Value = 1341.3334
$ws.Range("L46").Value = 3649.875
$ws.Range("M46").Value = -1153.3334
$ws.Range("N46").Value = -4025.875
$ws.Range("H55").Value = 553.45
$ws.Range("I55").Value = 519.63635
$ws.Range("K55").Value = 519.63635
$ws.Range("M55").Value = -346.63635

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H132").Value = 4049.5334
$ws.Range("J132").Value = 3093.3845
$ws.Range("L132").Value = 9280.1535
$ws.Range("N132").Value = -14340.1535
$ws.Range("H136").Value = 4567.1055
$ws.Range("I136").Value = 2628.8572
$ws.Range("J136").Value = 5697.75
$ws.Range("K136").Value = 7886.571599999999
$ws.Range("L136").Value = 17093.25
$ws.Range("M136").Value = -5336.571599999999
$ws.Range("N136").Value = -22193.25
